$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "ID", "`$1.00 Games", "Rollin' in the Dough", 16140, 1, "2019-03-12")
    ,@(3, "ID", "`$1.00 Games", "Lucky 13", 18197, 2, "2019-03-12")
    ,@(4, "ID", "`$1.00 Games", "Sheep Thrills", 19285, 2, "2019-03-12")
    ,@(5, "ID", "`$1.00 Games", "Rapid Refund", 18628, 1, "2019-03-12")
    ,@(6, "ID", "`$2.00 Games", "Lucky No 7 purple", 17251, 1, "2019-03-12")
    ,@(7, "ID", "`$2.00 Games", "Cherry, Lime, and Orange Twist", 19122, 2, "2019-03-12")
    ,@(8, "ID", "`$2.00 Games", "Lumberjack", 18629, 1, "2019-03-12")
    ,@(9, "ID", "`$2.00 Games", "Holiday Gold", 16141, 1, "2019-02-15")
    ,@(10, "ID", "`$3.00 Games", "Jolly Holiday Cashword", 16129, 1, "2019-03-12")
    ,@(11, "ID", "`$3.00 Games", "Slingo 3X", 18198, 2, "2019-03-12")
    ,@(12, "ID", "`$3.00 Games", "Bingo Plus", 15967, 1, "2019-03-12")
    ,@(13, "ID", "`$3.00 Games", "Chameleon Cashword", 19019, 2, "2019-03-12")
    ,@(14, "ID", "`$3.00 Games", "Bear Assets Bingo", 19286, 2, "2019-03-12")
    ,@(15, "ID", "`$5.00 Games", "Bengals", 2010, 3, "2019-03-12")
    ,@(16, "ID", "`$5.00 Games", "Winter Green", 16148, 1, "2019-03-07")
    ,@(17, "ID", "`$5.00 Games", "Marilyn Monroe", 18406, 1, "2019-03-12")
    ,@(18, "ID", "`$5.00 Games", "Plum Crazy", 15839, 1, "2019-03-12")
    ,@(19, "ID", "`$5.00 Games", "Crazy Cashword Connect", 19020, 2, "2019-03-12")
    ,@(20, "ID", "`$5.00 Games", "Triple Red 7s", 18405, 1, "2019-03-12")
    ,@(21, "ID", "`$5.00 Games", "Big Beach Party Bingo", 1778, 1, "2019-03-12")
    ,@(22, "ID", "`$5.00 Games", "Broncos", 2008, 2, "2019-03-12")
    ,@(23, "ID", "`$5.00 Games", "`$30,000 30th Anniversary", 19102, 2, "2019-03-12")
    ,@(24, "ID", "`$5.00 Games", "Vandals", 2009, 1, "2019-03-12")
    ,@(25, "ID", "`$5.00 Games", "5 Star Fortune", 16210, 0, "2019-03-04")
    ,@(26, "ID", "`$10.00 Games", "Tuxedo Cashword", 16036, 1, "2019-03-12")
    ,@(27, "ID", "`$10.00 Games", "`$100,000 30th Anniversary", 19144, 2, "2019-03-12")
    ,@(28, "ID", "`$10.00 Games", "Limited 4", 18396, 8595, "2019-03-12")
    ,@(29, "ID", "`$10.00 Games", "Lucky Jackpot Super Ticket", 15710, 1, "2019-03-12")
    ,@(30, "ID", "`$10.00 Games", "Limited 3rd Edition", 1786, 56, "2019-03-04")
    ,@(31, "ID", "`$20.00 Games", "Epic Fortune", 66, 1, "2019-02-12")
    ,@(32, "ID", "`$20.00 Games", "Juggernaut Jackpot", 15619, 1, "2019-03-12")
    ,@(33, "ID", "`$20.00 Games", "High Roller", 19021, 2, "2019-03-12")
    ,@(34, "ID", "`$20.00 Games", "Jackpot Frost", 16130, 0, "2019-03-05")
    ,@(35, "ID", "`$25.00 Games", "Tycoon Treasures", 16209, 1, "2019-03-12")
    ,@(36, "ID", "`$30.00 Games", "`$300000 30th Anniversary", 19145, 2, "2019-03-12")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $fcell = $ws.Cells.Item($r, 6)
    $fcell.NumberFormat = "@"
    $fcell.Value = $row[6]
    $fcell.ClearFormats()
}
